$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 831.3333
$ws.Range("I28").Value = 805.8421
$ws.Range("J28").Value = 855.55
$ws.Range("K28").Value = 805.8421
$ws.Range("L28").Value = 855.55
$ws.Range("M28").Value = -320.8421
$ws.Range("N28").Value = -1825.55
$ws.Range("H40").Value = 1575.25
$ws.Range("I40").Value = 1050
$ws.Range("J40").Value = 2100.5
$ws.Range("K40").Value = 1050
$ws.Range("L40").Value = 2100.5
$ws.Range("M40").Value = -875
$ws.Range("N40").Value = -2450.5
$ws.Range("H123").Value = 39765.453
$ws.Range("J123").Value = 39765.453
$ws.Range("L123").Value = 39765.453
$ws.Range("N123").Value = -49565.453
$ws.Range("H132").Value = 3581.9714
$ws.Range("I132").Value = 4380.3335
$ws.Range("J132").Value = 887.5
$ws.Range("K132").Value = 13141.0005
$ws.Range("L132").Value = 2662.5
$ws.Range("M132").Value = -10611.0005
$ws.Range("N132").Value = -7722.5
$ws.Range("H137").Value = 20001176
$ws.Range("I137").Value = 1092.8286
$ws.Range("J137").Value = 66668036
$ws.Range("K137").Value = 3278.4858
$ws.Range("L137").Value = 200004108
$ws.Range("M137").Value = -728.4858000000004
$ws.Range("N137").Value = -200009208
$ws.Range("H138").Value = 2889.5
$ws.Range("I138").Value = 2578.12
$ws.Range("J138").Value = 3070.535
$ws.Range("K138").Value = 7734.36
$ws.Range("L138").Value = 9211.605
$ws.Range("M138").Value = -2594.36
$ws.Range("N138").Value = -19491.605

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2247.05
$ws.Range("I32").Value = 1828.4337
$ws.Range("J32").Value = 4290.8823
$ws.Range("K32").Value = 1828.4337
$ws.Range("L32").Value = 4290.8823
$ws.Range("M32").Value = -1541.4337
$ws.Range("N32").Value = -4864.8823
$ws.Range("H45").Value = 1104.25
$ws.Range("I45").Value = 837.3333
$ws.Range("K45").Value = 837.3333
$ws.Range("M45").Value = -460.3333
$ws.Range("H61").Value = 2365225.8
$ws.Range("I61").Value = 2711155
$ws.Range("J61").Value = 1376
$ws.Range("K61").Value = 2711155
$ws.Range("L61").Value = 1376
$ws.Range("M61").Value = -2710943
$ws.Range("N61").Value = -1800
$ws.Range("H74").Value = 11909042
$ws.Range("I74").Value = 15625759
$ws.Range("J74").Value = 15550.2
$ws.Range("K74").Value = 15625759
$ws.Range("L74").Value = 15550.2
$ws.Range("M74").Value = -15624885
$ws.Range("N74").Value = -17298.2
$ws.Range("H77").Value = 11909042
$ws.Range("I77").Value = 15625759
$ws.Range("J77").Value = 15550.2
$ws.Range("K77").Value = 78128795
$ws.Range("L77").Value = 77751
$ws.Range("M77").Value = -78124427
$ws.Range("N77").Value = -86487
$ws.Range("H132").Value = 6565228
$ws.Range("I132").Value = 7555811
$ws.Range("K132").Value = 22667433
$ws.Range("M132").Value = -22664903
$ws.Range("H136").Value = 2365225.8
$ws.Range("I136").Value = 2711155
$ws.Range("J136").Value = 1376
$ws.Range("K136").Value = 8133465
$ws.Range("L136").Value = 4128
$ws.Range("M136").Value = -8130915
$ws.Range("N136").Value = -9228

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 32083
$ws.Range("J135").Value = 32083
$ws.Range("L135").Value = 32083
$ws.Range("N135").Value = -42223

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3183505.2
$ws.Range("I31").Value = 1078.4524
$ws.Range("J31").Value = 11537376
$ws.Range("K31").Value = 1078.4524
$ws.Range("L31").Value = 11537376
$ws.Range("M31").Value = -783.4523999999999
$ws.Range("N31").Value = -11537966
$ws.Range("H34").Value = 3183505.2
$ws.Range("I34").Value = 1078.4524
$ws.Range("J34").Value = 11537376
$ws.Range("K34").Value = 1078.4524
$ws.Range("L34").Value = 11537376
$ws.Range("M34").Value = -876.4523999999999
$ws.Range("N34").Value = -11537780
$ws.Range("H58").Value = 826.4177
$ws.Range("I58").Value = 794.24243
$ws.Range("J58").Value = 989.7692
$ws.Range("K58").Value = 794.24243
$ws.Range("L58").Value = 989.7692
$ws.Range("M58").Value = -591.24243
$ws.Range("N58").Value = -1395.7692
$ws.Range("H132").Value = 1775.625
$ws.Range("I132").Value = 1733.65
$ws.Range("J132").Value = 1985.5
$ws.Range("K132").Value = 5200.950000000001
$ws.Range("L132").Value = 5956.5
$ws.Range("M132").Value = -2670.950000000001
$ws.Range("N132").Value = -11016.5
$ws.Range("H134").Value = 1370.6389
$ws.Range("I134").Value = 1439.9
$ws.Range("J134").Value = 1024.3334
$ws.Range("K134").Value = 4319.700000000001
$ws.Range("L134").Value = 3073.0002
$ws.Range("M134").Value = -1784.700000000001
$ws.Range("N134").Value = -8143.0002
$ws.Range("H136").Value = 826.4177
$ws.Range("I136").Value = 794.24243
$ws.Range("J136").Value = 989.7692
$ws.Range("K136").Value = 2382.72729
$ws.Range("L136").Value = 2969.3076
$ws.Range("M136").Value = 167.2727100000002
$ws.Range("N136").Value = -8069.3076

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 45454784
$ws.Range("I33").Value = 50000230
$ws.Range("K33").Value = 300001380
$ws.Range("M33").Value = -300001097
$ws.Range("H34").Value = 1804.7333
$ws.Range("I34").Value = 143.54546
$ws.Range("J34").Value = 2766.4736
$ws.Range("K34").Value = 430.63638
$ws.Range("L34").Value = 8299.4208
$ws.Range("M34").Value = -346.63638
$ws.Range("N34").Value = -8467.4208
$ws.Range("H131").Value = 5727.1816
$ws.Range("I131").Value = 7078.8887
$ws.Range("J131").Value = 4791.385
$ws.Range("K131").Value = 21236.6661
$ws.Range("L131").Value = 14374.155
$ws.Range("M131").Value = -16196.6661
$ws.Range("N131").Value = -24454.155
$ws.Range("H132").Value = 100001650
$ws.Range("I132").Value = 200000800
$ws.Range("J132").Value = 2494
$ws.Range("K132").Value = 1800007200
$ws.Range("L132").Value = 22446
$ws.Range("M132").Value = -1800004670
$ws.Range("N132").Value = -27506
$ws.Range("H137").Value = 19291.635
$ws.Range("J137").Value = 25486.152
$ws.Range("L137").Value = 76458.45599999999
$ws.Range("N137").Value = -86658.45599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24392190
$ws.Range("I132").Value = 33335302
$ws.Range("J132").Value = 1884.3636
$ws.Range("K132").Value = 100005906
$ws.Range("L132").Value = 5653.0908
$ws.Range("M132").Value = -100003376
$ws.Range("N132").Value = -10713.0908
$ws.Range("H136").Value = 16712.875
$ws.Range("J136").Value = 16712.875
$ws.Range("L136").Value = 50138.625
$ws.Range("N136").Value = -55238.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H61").Value = 1807.7727
$ws.Range("I61").Value = 1709.25
$ws.Range("K61").Value = 1709.25
$ws.Range("M61").Value = -1507.25
$ws.Range("H113").Value = 1807.7727
$ws.Range("I113").Value = 1709.25
$ws.Range("K113").Value = 1709.25
$ws.Range("M113").Value = 460.75
$ws.Range("H132").Value = 3363.5833
$ws.Range("I132").Value = 3440.261
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 10320.783
$ws.Range("L132").Value = 4800
$ws.Range("M132").Value = -7790.782999999999
$ws.Range("N132").Value = -9860
$ws.Range("H136").Value = 1108.742
$ws.Range("I136").Value = 610.7778
$ws.Range("J136").Value = 4470
$ws.Range("K136").Value = 1832.3334
$ws.Range("L136").Value = 13410
$ws.Range("M136").Value = 717.6666
$ws.Range("N136").Value = -18510

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 40000000
$ws.Range("I4").Value = 10000000
$ws.Range("J4").Value = 100000000
$ws.Range("K4").Value = 10000000
$ws.Range("L4").Value = 100000000
$ws.Range("M4").Value = -9999887
$ws.Range("N4").Value = -100000226
$ws.Range("H54").Value = 18492
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 18492
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 18492
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = -19532
$ws.Range("H132").Value = 8347632.5
$ws.Range("I132").Value = 8718621
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 26155863
$ws.Range("L132").Value = 1200
$ws.Range("M132").Value = -26153333
$ws.Range("N132").Value = -6260
$ws.Range("H136").Value = 2806097.8
$ws.Range("I136").Value = 6874.909
$ws.Range("J136").Value = 7938006.5
$ws.Range("K136").Value = 20624.727
$ws.Range("L136").Value = 23814019.5
$ws.Range("M136").Value = -18074.727
$ws.Range("N136").Value = -23819119.5
